$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.045.36'
$ws.Range('E2').Value = '  -0.90%  '
$ws.Range('D3').Value = '3.419.66'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.73%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.88'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.50'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.86%  '
$ws.Range('D7').Value = '3.417.06'
$ws.Range('E7').Value = '  -0.88%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.467'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.135'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.78'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +5.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.404'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -4.17%  '
$ws.Range('D13').Value = '4.047.66'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000204'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.93'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.53%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.467.99'
$ws.Range('E16').Value = '  +0.80%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '66.759.58'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.116'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.49'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +6.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.09'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.82'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '416.12'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.70%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '78.08'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.581'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.92%  '
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').Value = '3.569.11'
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000111'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -6.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.20'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -5.66%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.45'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.63%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.76'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -6.53%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.160'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.22%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.47'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -7.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '24.59'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.70'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -4.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.52'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -9.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.61'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.01'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '172.37'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0860'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.07'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.872'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.90'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -12.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '45.76'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.69'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -8.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.16'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.06'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -5.15%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.27'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -6.30%  '
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.925'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.233'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.32%  '
